$wb = $excel.ActiveWorkbook

# --- Sheet ALC ---
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H11").Value = 46.375
$ws.Range("I11").Value = 46.375
$ws.Range("K11").Value = 46.375
$ws.Range("M11").Value = 93.625
$ws.Range("H21").Value = 11940.2
$ws.Range("I21").Value = 1500.5
$ws.Range("J21").Value = 18900
$ws.Range("K21").Value = 1500.5
$ws.Range("L21").Value = 18900
$ws.Range("M21").Value = -1032.5
$ws.Range("N21").Value = -19836
$ws.Range("H23").Value = 11940.2
$ws.Range("I23").Value = 1500.5
$ws.Range("J23").Value = 18900
$ws.Range("K23").Value = 1500.5
$ws.Range("L23").Value = 18900
$ws.Range("M23").Value = -1266.5
$ws.Range("N23").Value = -19368
$ws.Range("H29").Value = 403
$ws.Range("J29").Value = 0
$ws.Range("L29").Value = 0
$ws.Range("N29").Value = $null
$ws.Range("H38").Value = 386.33334
$ws.Range("I38").Value = 283.6
$ws.Range("K38").Value = 850.8000000000001
$ws.Range("M38").Value = -478.8000000000001
$ws.Range("H51").Value = 2993.3333
$ws.Range("I51").Value = 4980
$ws.Range("J51").Value = 2000
$ws.Range("K51").Value = 4980
$ws.Range("L51").Value = 2000
$ws.Range("M51").Value = -4496
$ws.Range("N51").Value = -2968
$ws.Range("H53").Value = 219.18182
$ws.Range("I53").Value = 238.5
$ws.Range("K53").Value = 238.5
$ws.Range("M53").Value = 398.5
$ws.Range("H58").Value = 843.7646999999999
$ws.Range("I58").Value = 843.7646999999999
$ws.Range("J58").Value = 0
$ws.Range("K58").Value = 2531.2941
$ws.Range("L58").Value = 0
$ws.Range("M58").Value = -2381.2941
$ws.Range("N58").Value = $null
$ws.Range("H86").Value = 107144840
$ws.Range("I86").Value = 250001400
$ws.Range("J86").Value = 2424.25
$ws.Range("K86").Value = 250001400
$ws.Range("L86").Value = 2424.25
$ws.Range("M86").Value = -250000277
$ws.Range("N86").Value = -4670.25
$ws.Range("H89").Value = 107144840
$ws.Range("I89").Value = 250001400
$ws.Range("J89").Value = 2424.25
$ws.Range("K89").Value = 1250007000
$ws.Range("L89").Value = 12121.25
$ws.Range("M89").Value = -1250001384
$ws.Range("N89").Value = -23353.25
$ws.Range("H132").Value = 3688.9546
$ws.Range("I132").Value = 3192.5789
$ws.Range("J132").Value = 6832.6665
$ws.Range("K132").Value = 9577.736699999999
$ws.Range("L132").Value = 20497.9995
$ws.Range("M132").Value = -7047.736699999999
$ws.Range("N132").Value = -25557.9995
$ws.Range("H138").Value = 2115.869
$ws.Range("I138").Value = 1695.7693
$ws.Range("J138").Value = 2427.9429
$ws.Range("K138").Value = 5087.3079
$ws.Range("L138").Value = 7283.8287
$ws.Range("M138").Value = 52.69210000000021
$ws.Range("N138").Value = -17563.8287
$ws.Range("H141").Value = 3767.4285
$ws.Range("I141").Value = 1177.4546
$ws.Range("J141").Value = 13264
$ws.Range("K141").Value = 3532.3638
$ws.Range("L141").Value = 39792
$ws.Range("M141").Value = 1647.6362
$ws.Range("N141").Value = -50152

# --- Sheet ARM ---
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H9").Value = 5000
$ws.Range("I9").Value = 0
$ws.Range("K9").Value = 0
$ws.Range("M9").Value = $null
$ws.Range("H20").Value = 5000
$ws.Range("I20").Value = 0
$ws.Range("K20").Value = 0
$ws.Range("M20").Value = $null
$ws.Range("H23").Value = 12547.81
$ws.Range("J23").Value = 9674.9
$ws.Range("L23").Value = 9674.9
$ws.Range("N23").Value = -10192.9
$ws.Range("H37").Value = 19999
$ws.Range("I37").Value = 0
$ws.Range("J37").Value = 19999
$ws.Range("K37").Value = 0
$ws.Range("L37").Value = 19999
$ws.Range("M37").Value = $null
$ws.Range("N37").Value = -20545
$ws.Range("H44").Value = 12569.857
$ws.Range("J44").Value = 12569.857
$ws.Range("L44").Value = 12569.857
$ws.Range("N44").Value = -13545.857
$ws.Range("H45").Value = 4324.8
$ws.Range("I45").Value = 4409.9
$ws.Range("K45").Value = 4409.9
$ws.Range("M45").Value = -4032.9

# --- Sheet BSM ---
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H35").Value = 17857
$ws.Range("J35").Value = 17857
$ws.Range("L35").Value = 17857
$ws.Range("N35").Value = -18477

# --- Sheet CRP ---
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H17").Value = 3000
$ws.Range("I17").Value = 3000
$ws.Range("K17").Value = 3000
$ws.Range("M17").Value = -2826
$ws.Range("H25").Value = 3500
$ws.Range("I25").Value = 3500
$ws.Range("K25").Value = 3500
$ws.Range("M25").Value = -3326
$ws.Range("H41").Value = 15499.25
$ws.Range("I41").Value = 2000
$ws.Range("K41").Value = 2000
$ws.Range("M41").Value = -1572
$ws.Range("H59").Value = 13623.125
$ws.Range("I59").Value = 3000
$ws.Range("J59").Value = 15140.714
$ws.Range("K59").Value = 3000
$ws.Range("L59").Value = 15140.714
$ws.Range("M59").Value = -1855
$ws.Range("N59").Value = -17430.714
$ws.Range("H60").Value = 9650.799999999999
$ws.Range("I60").Value = 4500
$ws.Range("J60").Value = 10223.111
$ws.Range("K60").Value = 4500
$ws.Range("L60").Value = 10223.111
$ws.Range("M60").Value = -3989
$ws.Range("N60").Value = -11245.111

# --- Sheet CUL ---
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H4").Value = 20003000
$ws.Range("I4").Value = 28003000
$ws.Range("J4").Value = 3000
$ws.Range("K4").Value = 84009000
$ws.Range("L4").Value = 9000
$ws.Range("M4").Value = -84008888
$ws.Range("N4").Value = -9224
$ws.Range("H107").Value = 587.6667
$ws.Range("J107").Value = 859.8
$ws.Range("L107").Value = 2579.4
$ws.Range("N107").Value = -6419.4

# --- Sheet GSM ---
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H5").Value = 4250
$ws.Range("I5").Value = 2000
$ws.Range("K5").Value = 2000
$ws.Range("M5").Value = -1888
$ws.Range("H43").Value = 4726.2144
$ws.Range("I43").Value = 742.7273
$ws.Range("K43").Value = 742.7273
$ws.Range("M43").Value = -591.7273
$ws.Range("H46").Value = 4040.8696
$ws.Range("J46").Value = 4349.95
$ws.Range("L46").Value = 4349.95
$ws.Range("N46").Value = -4661.95
$ws.Range("H57").Value = 16714.143
$ws.Range("J57").Value = 16714.143
$ws.Range("L57").Value = 16714.143
$ws.Range("N57").Value = -18354.143
$ws.Range("H93").Value = 0
$ws.Range("J93").Value = 0
$ws.Range("L93").Value = 0
$ws.Range("N93").Value = $null

# --- Sheet LTW ---
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H2").Value = 47334.332
$ws.Range("I2").Value = 1999
$ws.Range("J2").Value = 70002
$ws.Range("K2").Value = 1999
$ws.Range("L2").Value = 70002
$ws.Range("M2").Value = -1887
$ws.Range("N2").Value = -70226
$ws.Range("H20").Value = 61603.6
$ws.Range("J20").Value = 61603.6
$ws.Range("L20").Value = 61603.6
$ws.Range("N20").Value = -62055.6
$ws.Range("H46").Value = 900
$ws.Range("J46").Value = 920
$ws.Range("L46").Value = 920
$ws.Range("N46").Value = -1296
